$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179
$ws.Range("J2").Value = 786
$ws.Range("L2").Value = 786
$ws.Range("N2").Value = -1012
$ws.Range("H4").Value = 7608.2856
$ws.Range("I4").Value = 4314.75
$ws.Range("K4").Value = 4314.75
$ws.Range("M4").Value = -4200.75
$ws.Range("H19").Value = 16454.111
$ws.Range("I19").Value = 2384.1428
$ws.Range("J19").Value = 25407.727
$ws.Range("K19").Value = 2384.1428
$ws.Range("L19").Value = 25407.727
$ws.Range("M19").Value = -2209.1428
$ws.Range("N19").Value = -25757.727
$ws.Range("H32").Value = 8587.091
$ws.Range("I32").Value = 8660
$ws.Range("J32").Value = 8559.75
$ws.Range("K32").Value = 8660
$ws.Range("L32").Value = 8559.75
$ws.Range("M32").Value = -8334
$ws.Range("N32").Value = -9211.75
$ws.Range("H40").Value = 2742.5715
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H53").Value = 66937.47
$ws.Range("J53").Value = 200606
$ws.Range("L53").Value = 200606
$ws.Range("N53").Value = -201880
$ws.Range("H64").Value = 4107.4287
$ws.Range("I64").Value = 3501
$ws.Range("J64").Value = 4350
$ws.Range("K64").Value = 3501
$ws.Range("L64").Value = 4350
$ws.Range("M64").Value = -3253
$ws.Range("N64").Value = -4846
$ws.Range("H67").Value = 4107.4287
$ws.Range("I67").Value = 3501
$ws.Range("J67").Value = 4350
$ws.Range("K67").Value = 3501
$ws.Range("L67").Value = 4350
$ws.Range("M67").Value = -2643
$ws.Range("N67").Value = -6066
$ws.Range("H70").Value = 602535.4
$ws.Range("I70").Value = 852573.4399999999
$ws.Range("K70").Value = 2557720.32
$ws.Range("M70").Value = -2557450.32
$ws.Range("H73").Value = 602535.4
$ws.Range("I73").Value = 852573.4399999999
$ws.Range("K73").Value = 2557720.32
$ws.Range("M73").Value = -2556784.32
$ws.Range("H76").Value = 5274861.5
$ws.Range("I76").Value = 2689.8
$ws.Range("J76").Value = 7671303.5
$ws.Range("K76").Value = 2689.8
$ws.Range("L76").Value = 7671303.5
$ws.Range("M76").Value = -2374.8
$ws.Range("N76").Value = -7671933.5
$ws.Range("H79").Value = 5274861.5
$ws.Range("I79").Value = 2689.8
$ws.Range("J79").Value = 7671303.5
$ws.Range("K79").Value = 2689.8
$ws.Range("L79").Value = 7671303.5
$ws.Range("M79").Value = -1597.8
$ws.Range("N79").Value = -7673487.5
$ws.Range("H88").Value = 716
$ws.Range("I88").Value = 847.6
$ws.Range("J88").Value = 633.75
$ws.Range("K88").Value = 847.6
$ws.Range("L88").Value = 633.75
$ws.Range("M88").Value = -441.6
$ws.Range("N88").Value = -1445.75
$ws.Range("H91").Value = 716
$ws.Range("I91").Value = 847.6
$ws.Range("J91").Value = 633.75
$ws.Range("K91").Value = 847.6
$ws.Range("L91").Value = 633.75
$ws.Range("M91").Value = 556.4
$ws.Range("N91").Value = -3441.75
$ws.Range("H98").Value = 2185.4
$ws.Range("I98").Value = 2402.6155
$ws.Range("K98").Value = 2402.6155
$ws.Range("M98").Value = -904.6154999999999
$ws.Range("H122").Value = 2185.4
$ws.Range("I122").Value = 2402.6155
$ws.Range("K122").Value = 7207.8465
$ws.Range("M122").Value = -4757.8465
$ws.Range("H131").Value = 9714.615
$ws.Range("I131").Value = 3063.3333
$ws.Range("J131").Value = 11710
$ws.Range("K131").Value = 9189.999899999999
$ws.Range("L131").Value = 35130
$ws.Range("M131").Value = -4149.999899999999
$ws.Range("N131").Value = -45210
$ws.Range("H135").Value = 2249.25
$ws.Range("I135").Value = 999
$ws.Range("J135").Value = 3499.5
$ws.Range("K135").Value = 8991
$ws.Range("L135").Value = 31495.5
$ws.Range("M135").Value = -6456
$ws.Range("N135").Value = -36565.5
$ws.Range("H137").Value = 2438.2942
$ws.Range("I137").Value = 1970.44
$ws.Range("J137").Value = 3737.889
$ws.Range("K137").Value = 5911.32
$ws.Range("L137").Value = 11213.667
$ws.Range("M137").Value = -3361.32
$ws.Range("N137").Value = -16313.667
$ws.Range("H138").Value = 1943.1
$ws.Range("I138").Value = 1464.1111
$ws.Range("K138").Value = 4392.3333
$ws.Range("M138").Value = 747.6666999999998
$ws.Range("H141").Value = 6471
$ws.Range("I141").Value = 6471
$ws.Range("K141").Value = 19413
$ws.Range("M141").Value = -14233
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 811.5789
$ws.Range("I2").Value = 898.125
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 898.125
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -785.125
$ws.Range("N2").Value = -576
$ws.Range("H19").Value = 3700
$ws.Range("I19").Value = 3700
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 3700
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -3471
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 6900.1313
$ws.Range("I32").Value = 6658.9707
$ws.Range("J32").Value = 8950
$ws.Range("K32").Value = 6658.9707
$ws.Range("L32").Value = 8950
$ws.Range("M32").Value = -6371.9707
$ws.Range("N32").Value = -9524
$ws.Range("H63").Value = 66671176
$ws.Range("I63").Value = 111113300
$ws.Range("J63").Value = 22229056
$ws.Range("K63").Value = 111113300
$ws.Range("L63").Value = 22229056
$ws.Range("M63").Value = -111112614
$ws.Range("N63").Value = -22230428
$ws.Range("H66").Value = 66671176
$ws.Range("I66").Value = 111113300
$ws.Range("J66").Value = 22229056
$ws.Range("K66").Value = 555566500
$ws.Range("L66").Value = 111145280
$ws.Range("M66").Value = -555563068
$ws.Range("N66").Value = -111152144
$ws.Range("H74").Value = 90913230
$ws.Range("I74").Value = 100003656
$ws.Range("K74").Value = 100003656
$ws.Range("M74").Value = -100002782
$ws.Range("H77").Value = 90913230
$ws.Range("I77").Value = 100003656
$ws.Range("K77").Value = 500018280
$ws.Range("M77").Value = -500013912
$ws.Range("H80").Value = 69788.664
$ws.Range("I80").Value = 69788.664
$ws.Range("K80").Value = 69788.664
$ws.Range("M80").Value = -68790.664
$ws.Range("H83").Value = 69788.664
$ws.Range("I83").Value = 69788.664
$ws.Range("K83").Value = 209365.992
$ws.Range("M83").Value = -204373.992
$ws.Range("H88").Value = 16668112
$ws.Range("I88").Value = 41667580
$ws.Range("J88").Value = 1799.5
$ws.Range("K88").Value = 41667580
$ws.Range("L88").Value = 1799.5
$ws.Range("M88").Value = -41667174
$ws.Range("N88").Value = -2611.5
$ws.Range("H91").Value = 16668112
$ws.Range("I91").Value = 41667580
$ws.Range("J91").Value = 1799.5
$ws.Range("K91").Value = 41667580
$ws.Range("L91").Value = 1799.5
$ws.Range("M91").Value = -41666176
$ws.Range("N91").Value = -4607.5
$ws.Range("H102").Value = 3031383.8
$ws.Range("I102").Value = 3717496.2
$ws.Range("K102").Value = 3717496.2
$ws.Range("M102").Value = -3715874.2
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H116").Value = 811.5789
$ws.Range("I116").Value = 898.125
$ws.Range("J116").Value = 350
$ws.Range("K116").Value = 898.125
$ws.Range("L116").Value = 350
$ws.Range("M116").Value = 1395.875
$ws.Range("N116").Value = -4938
$ws.Range("H122").Value = 1921.18
$ws.Range("I122").Value = 1161
$ws.Range("K122").Value = 3483
$ws.Range("M122").Value = -1033
$ws.Range("H132").Value = 3310.2917
$ws.Range("I132").Value = 1803.9286
$ws.Range("J132").Value = 5419.2
$ws.Range("K132").Value = 5411.7858
$ws.Range("L132").Value = 16257.6
$ws.Range("M132").Value = -2881.7858
$ws.Range("N132").Value = -21317.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 811.5789
$ws.Range("I3").Value = 898.125
$ws.Range("J3").Value = 350
$ws.Range("K3").Value = 898.125
$ws.Range("L3").Value = 350
$ws.Range("M3").Value = -784.125
$ws.Range("N3").Value = -578
$ws.Range("H20").Value = 18254.334
$ws.Range("I20").Value = 20678.076
$ws.Range("K20").Value = 20678.076
$ws.Range("M20").Value = -20431.076
$ws.Range("H82").Value = 4036.5715
$ws.Range("I82").Value = 4036.5715
$ws.Range("K82").Value = 4036.5715
$ws.Range("M82").Value = -3653.5715
$ws.Range("H85").Value = 4036.5715
$ws.Range("I85").Value = 4036.5715
$ws.Range("K85").Value = 4036.5715
$ws.Range("M85").Value = -2710.5715
$ws.Range("H86").Value = 1719.5
$ws.Range("I86").Value = 1243
$ws.Range("K86").Value = 1243
$ws.Range("M86").Value = -120
$ws.Range("H89").Value = 1719.5
$ws.Range("I89").Value = 1243
$ws.Range("K89").Value = 6215
$ws.Range("M89").Value = -599
$ws.Range("H94").Value = 20836096
$ws.Range("I94").Value = 27779794
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 27779794
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -27779343
$ws.Range("N94").Value = -5902
$ws.Range("H105").Value = 2169.1072
$ws.Range("I105").Value = 1935.6666
$ws.Range("K105").Value = 1935.6666
$ws.Range("M105").Value = -188.6666
$ws.Range("H107").Value = 41836290
$ws.Range("I107").Value = 334748.16
$ws.Range("J107").Value = 83337830
$ws.Range("K107").Value = 334748.16
$ws.Range("L107").Value = 83337830
$ws.Range("M107").Value = -332828.16
$ws.Range("N107").Value = -83341670
$ws.Range("H134").Value = 4728.727
$ws.Range("I134").Value = 4358.4443
$ws.Range("K134").Value = 13075.3329
$ws.Range("M134").Value = -10540.3329
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 844
$ws.Range("I22").Value = 613
$ws.Range("K22").Value = 613
$ws.Range("M22").Value = -263
$ws.Range("H31").Value = 3460.7856
$ws.Range("I31").Value = 3336
$ws.Range("J31").Value = 3627.1667
$ws.Range("K31").Value = 3336
$ws.Range("L31").Value = 3627.1667
$ws.Range("M31").Value = -3041
$ws.Range("N31").Value = -4217.1667
$ws.Range("H34").Value = 3460.7856
$ws.Range("I34").Value = 3336
$ws.Range("J34").Value = 3627.1667
$ws.Range("K34").Value = 3336
$ws.Range("L34").Value = 3627.1667
$ws.Range("M34").Value = -3134
$ws.Range("N34").Value = -4031.1667
$ws.Range("H58").Value = 4442.385
$ws.Range("I58").Value = 2153.2856
$ws.Range("J58").Value = 7113
$ws.Range("K58").Value = 2153.2856
$ws.Range("L58").Value = 7113
$ws.Range("M58").Value = -1950.2856
$ws.Range("N58").Value = -7519
$ws.Range("H62").Value = 139448.5
$ws.Range("I62").Value = 3897
$ws.Range("J62").Value = 275000
$ws.Range("K62").Value = 3897
$ws.Range("L62").Value = 275000
$ws.Range("M62").Value = -3273
$ws.Range("N62").Value = -276248
$ws.Range("H65").Value = 139448.5
$ws.Range("I65").Value = 3897
$ws.Range("J65").Value = 275000
$ws.Range("K65").Value = 19485
$ws.Range("L65").Value = 1375000
$ws.Range("M65").Value = -16365
$ws.Range("N65").Value = -1381240
$ws.Range("H107").Value = 460.96875
$ws.Range("I107").Value = 439.1111
$ws.Range("K107").Value = 439.1111
$ws.Range("M107").Value = 1480.8889
$ws.Range("H132").Value = 4131.778
$ws.Range("I132").Value = 2716
$ws.Range("K132").Value = 8148
$ws.Range("M132").Value = -5618
$ws.Range("H136").Value = 4442.385
$ws.Range("I136").Value = 2153.2856
$ws.Range("J136").Value = 7113
$ws.Range("K136").Value = 6459.8568
$ws.Range("L136").Value = 21339
$ws.Range("M136").Value = -3909.8568
$ws.Range("N136").Value = -26439
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1327.6875
$ws.Range("J2").Value = 1101.2
$ws.Range("L2").Value = 6607.200000000001
$ws.Range("N2").Value = -6833.200000000001
$ws.Range("H5").Value = 2310.923
$ws.Range("I5").Value = 436.42856
$ws.Range("K5").Value = 1309.28568
$ws.Range("M5").Value = -1197.28568
$ws.Range("H34").Value = 749.2222
$ws.Range("I34").Value = 347
$ws.Range("J34").Value = 864.1429000000001
$ws.Range("K34").Value = 1041
$ws.Range("L34").Value = 2592.4287
$ws.Range("M34").Value = -957
$ws.Range("N34").Value = -2760.4287
$ws.Range("H36").Value = 333933.34
$ws.Range("I36").Value = 333933.34
$ws.Range("K36").Value = 1001800.02
$ws.Range("M36").Value = -1001631.02
$ws.Range("H55").Value = 2875
$ws.Range("J55").Value = 3500
$ws.Range("L55").Value = 10500
$ws.Range("N55").Value = -10854
$ws.Range("H131").Value = 3054.247
$ws.Range("I131").Value = 1938.6
$ws.Range("J131").Value = 3123.975
$ws.Range("K131").Value = 5815.799999999999
$ws.Range("L131").Value = 9371.924999999999
$ws.Range("M131").Value = -775.7999999999993
$ws.Range("N131").Value = -19451.925
$ws.Range("H135").Value = 2310.923
$ws.Range("I135").Value = 436.42856
$ws.Range("K135").Value = 3927.85704
$ws.Range("M135").Value = -1392.85704
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8629.053
$ws.Range("I70").Value = 7330.75
$ws.Range("J70").Value = 10854.714
$ws.Range("K70").Value = 7330.75
$ws.Range("L70").Value = 10854.714
$ws.Range("M70").Value = -7060.75
$ws.Range("N70").Value = -11394.714
$ws.Range("H73").Value = 8629.053
$ws.Range("I73").Value = 7330.75
$ws.Range("J73").Value = 10854.714
$ws.Range("K73").Value = 7330.75
$ws.Range("L73").Value = 10854.714
$ws.Range("M73").Value = -6394.75
$ws.Range("N73").Value = -12726.714
$ws.Range("H132").Value = 3575.7104
$ws.Range("I132").Value = 3407.2334
$ws.Range("J132").Value = 4207.5
$ws.Range("K132").Value = 10221.7002
$ws.Range("L132").Value = 12622.5
$ws.Range("M132").Value = -7691.700199999999
$ws.Range("N132").Value = -17682.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19234098
$ws.Range("J7").Value = 3937
$ws.Range("L7").Value = 3937
$ws.Range("N7").Value = -4161
$ws.Range("H16").Value = 456.5
$ws.Range("I16").Value = 219.375
$ws.Range("K16").Value = 219.375
$ws.Range("M16").Value = -49.375
$ws.Range("H22").Value = 2600.5
$ws.Range("I22").Value = 2600.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2600.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2305.5
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 2600.5
$ws.Range("I27").Value = 2600.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2600.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2493.5
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 2340
$ws.Range("J46").Value = 2340
$ws.Range("L46").Value = 2340
$ws.Range("N46").Value = -2716
$ws.Range("H68").Value = 4896.375
$ws.Range("I68").Value = 4116.25
$ws.Range("J68").Value = 5676.5
$ws.Range("K68").Value = 4116.25
$ws.Range("L68").Value = 5676.5
$ws.Range("M68").Value = -3367.25
$ws.Range("N68").Value = -7174.5
$ws.Range("H71").Value = 4896.375
$ws.Range("I71").Value = 4116.25
$ws.Range("J71").Value = 5676.5
$ws.Range("K71").Value = 20581.25
$ws.Range("L71").Value = 28382.5
$ws.Range("M71").Value = -16837.25
$ws.Range("N71").Value = -35870.5
$ws.Range("H74").Value = 49950
$ws.Range("I74").Value = 49950
$ws.Range("K74").Value = 49950
$ws.Range("M74").Value = -48952
$ws.Range("H77").Value = 49950
$ws.Range("I77").Value = 49950
$ws.Range("K77").Value = 149850
$ws.Range("M77").Value = -144858
$ws.Range("H93").Value = 25642988
$ws.Range("I93").Value = 83334840
$ws.Range("K93").Value = 83334840
$ws.Range("M93").Value = -83333592
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 3857.6052
$ws.Range("I122").Value = 3628.348
$ws.Range("K122").Value = 10885.044
$ws.Range("M122").Value = -8435.044
$ws.Range("H126").Value = 19234098
$ws.Range("J126").Value = 3937
$ws.Range("L126").Value = 11811
$ws.Range("N126").Value = -16751
$ws.Range("H132").Value = 20004408
$ws.Range("I132").Value = 25643496
$ws.Range("K132").Value = 76930488
$ws.Range("M132").Value = -76927958
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 6680666.5
$ws.Range("J48").Value = 21000
$ws.Range("L48").Value = 21000
$ws.Range("N48").Value = -22138
$ws.Range("H94").Value = 16250
$ws.Range("J94").Value = 16250
$ws.Range("L94").Value = 16250
$ws.Range("N94").Value = -18052
$ws.Range("H114").Value = 94999
$ws.Range("J114").Value = 94999
$ws.Range("L114").Value = 94999
$ws.Range("N114").Value = -103677
$ws.Range("H122").Value = 1809.0571
$ws.Range("I122").Value = 1577.6538
$ws.Range("K122").Value = 4732.9614
$ws.Range("M122").Value = -2282.9614
$ws.Range("H123").Value = 42198
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 5152.8647
$ws.Range("I132").Value = 4898.6562
$ws.Range("K132").Value = 14695.9686
$ws.Range("M132").Value = -12165.9686
